$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1/2, merged like the other header columns) ---
$ws.Range("AK1").Value = "GridLayers"
$ws.Range("AL1").Value = "GridStyleURL"

$ws.Range("AK1:AK2").Merge()
$ws.Range("AL1:AL2").Merge()

# Match the formatting of the existing header cells (center/center, same font)
$ws.Range("AH1").Copy()
$ws.Range("AK1:AL2").PasteSpecial(-4122)

# --- New data columns: GridLayers (AK) / GridStyleURL (AL) ---
$ws.Range("AK3").Value = "montreal_grid"
$ws.Range("AL3").Value = "cl1qihonq002415ng3q7njruw"

$ws.Range("AK4").Value = "vienna_grid"
$ws.Range("AL4").Value = "cl1qj5iio000h14lsrpfbbbs7"

$ws.Range("AK5").Value = "barcelona_grid"
$ws.Range("AL5").Value = "cl1qjcbr9000e15s7ids0qcrt"

$ws.Range("AK6").Value = "budapest_grid"
$ws.Range("AL6").Value = "cl1qjqiij002h14s6sigeccvb"

$ws.Range("AK7").Value = "quebec_grid"
$ws.Range("AL7").Value = "cl1qjjphy000g15s7jf1t62sy"

# Match the formatting of the existing data cells in that row band (center/center)
$ws.Range("AH3").Copy()
$ws.Range("AK3:AL7").PasteSpecial(-4122)

# Size the two new columns similarly to their neighbours
$ws.Columns("AK:AL").AutoFit()

$ws.Application.CutCopyMode = $false

$ws.Range("A1").Select()
